$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, copying header style from E1
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Populate time_taken values for each data row (F2:F30)
$ws.Range("F2").Value = "2021-10-05 13:40:00.905891"
$ws.Range("F3").Value = "2021-10-05 13:40:00.905904"
$ws.Range("F4").Value = "2021-10-05 13:40:00.905908"
$ws.Range("F5").Value = "2021-10-05 13:40:00.905912"
$ws.Range("F6").Value = "2021-10-05 13:40:00.905915"
$ws.Range("F7").Value = "2021-10-05 13:40:00.905918"
$ws.Range("F8").Value = "2021-10-05 13:40:00.905921"
$ws.Range("F9").Value = "2021-10-05 13:40:00.905924"
$ws.Range("F10").Value = "2021-10-05 13:40:00.905928"
$ws.Range("F11").Value = "2021-10-05 13:40:00.905931"
$ws.Range("F12").Value = "2021-10-05 13:40:00.905934"
$ws.Range("F13").Value = "2021-10-05 13:40:00.905937"
$ws.Range("F14").Value = "2021-10-05 13:40:00.905940"
$ws.Range("F15").Value = "2021-10-05 13:40:00.905943"
$ws.Range("F16").Value = "2021-10-05 13:40:00.905946"
$ws.Range("F17").Value = "2021-10-05 13:40:00.905949"
$ws.Range("F18").Value = "2021-10-05 13:40:00.905953"
$ws.Range("F19").Value = "2021-10-05 13:40:00.905956"
$ws.Range("F20").Value = "2021-10-05 13:40:00.905959"
$ws.Range("F21").Value = "2021-10-05 13:40:00.905962"
$ws.Range("F22").Value = "2021-10-05 13:40:00.905965"
$ws.Range("F23").Value = "2021-10-05 13:40:00.905968"
$ws.Range("F24").Value = "2021-10-05 13:40:00.905971"
$ws.Range("F25").Value = "2021-10-05 13:40:00.905974"
$ws.Range("F26").Value = "2021-10-05 13:40:00.905978"
$ws.Range("F27").Value = "2021-10-05 13:40:00.905981"
$ws.Range("F28").Value = "2021-10-05 13:40:00.905984"
$ws.Range("F29").Value = "2021-10-05 13:40:00.905992"
$ws.Range("F30").Value = "2021-10-05 13:40:00.905996"

$excel.CutCopyMode = 0
